$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I28").Value = 3
